$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.287.06"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "3.149.50"
$ws.Range("E3").Value = "  +2.77%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.43"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.96"
$ws.Range("E6").Value = "  +5.37%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.140.84"
$ws.Range("E8").Value = "  +2.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  +16.59%  "

$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  +1.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.49"
$ws.Range("E13").Value = "  +2.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000227"
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("D15").Value = "3.652.67"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "65.277.79"
$ws.Range("E16").Value = "  +1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "539.06"
$ws.Range("E17").Value = "  +10.41%  "

$ws.Range("D19").Value = "3.149.15"
$ws.Range("E19").Value = "  +2.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.85"
$ws.Range("E20").Value = "  +3.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  +2.48%  "

$ws.Range("E22").Value = "  +4.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +4.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.05"
$ws.Range("E24").Value = "  +4.35%  "

$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  +16.07%  "

$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  +5.01%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  +3.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.56"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  +3.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "557.20"
$ws.Range("E34").Value = "  +8.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.52"
$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.18"
$ws.Range("E36").Value = "  +4.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0450"
$ws.Range("E37").Value = "  +10.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.17"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0833"
$ws.Range("E39").Value = "  +5.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +9.26%  "

$ws.Range("E41").Value = "  +2.92%  "

$ws.Range("D42").Value = "3.070.83"
$ws.Range("E42").Value = "  +6.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.38"
$ws.Range("E43").Value = "  +0.80%  "

$ws.Range("E44").Value = "  +7.72%  "

$ws.Range("E45").Value = "  +8.26%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.46"
$ws.Range("E47").Value = "  +2.76%  "

$ws.Range("E48").Value = "  -2.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.111"
$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "119.89"
$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  +4.19%  "
